# Apply the "members" sheet re-sort / data-cleanup edit described in the
# commit:
#   - remove the "emmy thornton" and "joy maina" rows (no longer members)
#   - update a couple of photo filenames from HEIC/JPEG to PNG
#   - add a profile photo for ellie mattox
#   - backfill missing bios for several undergrads
# The director sheet's C2 cell keeps the exact same string value; the
# underlying diff there is only a shared-string-table renumbering side
# effect of removing the two rows' worth of now-unused strings, so no
# explicit edit is required there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("members")

# --- Remove the two members who left (delete from the bottom up so the
#     earlier row number, 19, isn't invalidated by the later delete). ---
$ws.Rows.Item(22).Delete()   # "joy maina"
$ws.Rows.Item(19).Delete()   # "emmy thornton"

# --- Photo filename fixes (HEIC/JPEG -> PNG) and new photo for ellie ---
# (written ellie, then ben, then claudia, to match the order these new
#  filenames were introduced into the workbook's shared-string table)
$ws.Range("A18").Value = "ellie_mattox.png"
$ws.Range("A17").Value = "ben_hanson.png"
$ws.Range("A4").Value = "claudia_gonciulea.png"

# --- Backfill missing bios ---
$ws.Range("E5").Value = "Jennifer is planning to double major in Quantitative Social Science and Comparative Literature. She has a background in mathematics and programming, and is interested in the humanities and literature. She hopes to use machine learning to enhance our understanding of humanitarian issues."
$ws.Range("E7").Value = "Sarah is interested in cognitive science, psychology, and neuroscience. She is working on using EEG to track how people learn from online course videos."
$ws.Range("E8").Value = "Chelsea is a Computer Science and Neuroscience double major. She is interested in artificial intelligence, machine learning, and building computational models."
$ws.Range("E9").Value = "Jacob's work is focused on using machine learning and deep learning models to understand mental health. He is also interested in how cognitive and mental function relate to one another."
$ws.Range("E10").Value = 'Aidan is a pre-med student interested in how we can learn as quickly and efficiently as possible. He is a dedicated practitioner of memory "hacks" like the chaine method, story method, and method of loci.'
$ws.Range("E11").Value = "Alexandra is interested in education technology and brings to the lab her extensive computer science experience. She is especially interested in helping students with learning differences, and in developing brain-based learning tools."
$ws.Range("E16").Value = "Azaire is a Mathematics and Anthropology double major. She is particularly interested in understanding complex systems like brain networks and financial markets."
$ws.Range("E19").Value = "Evan is a Mathematics major and Computer Science minor. He is interested in algorithmic trading, algorithmic problem solving, and building models of financial market dynamics."
$ws.Range("E20").Value = "Jax is a Mathematics and Economics double major. He is interested in using EEG to understand how people learn new concepts and skills."
$ws.Range("E22").Value = "Luca is a Neuroscience major who is passionate about how we can use technology and neuroscience to enhance learning in everyday life. He is also interested in wearable brain recording devices."
$ws.Range("E24").Value = "Sam is a Quantitative Social Science and Economics double major. He is interested in how social peer influences shape beliefs and affect how people make decisions."

# --- Row heights grow now that those rows carry bio text ---
$ws.Rows.Item(5).RowHeight = 96
$ws.Rows.Item(7).RowHeight = 48
$ws.Rows.Item(8).RowHeight = 64
$ws.Rows.Item(9).RowHeight = 64
$ws.Rows.Item(10).RowHeight = 64
$ws.Rows.Item(11).RowHeight = 80
$ws.Rows.Item(16).RowHeight = 64
$ws.Rows.Item(19).RowHeight = 64
$ws.Rows.Item(20).RowHeight = 48
$ws.Rows.Item(22).RowHeight = 80
$ws.Rows.Item(24).RowHeight = 64

# --- Selection ends up on the "evan mcdermid" row (whole row selected,
#     with no pinned top-left cell) after the edit ---
$ws.Range("A19:F19").EntireRow.Select()
